$data = @(
    @("U mnie to nie zadziałało.", "polish"),
    @("Czy moge prosić o jakąś podpowiedź?", "polish"),
    @("Teraz możesz zobaczyć się ze swoim lekarzem bez wychodzenia z domu.", "polish"),
    @("Codziennie po lekcjach gram w tenisa ze swoim najlepszym kolegą.", "polish"),
    @("Prawdziwych przyjaciół poznajemy w biedzie.", "polish"),
    @("Dzien dobry, proszę o informację, jak otworzyć plik na dysku.", "polish"),
    @("Każdego dnia ćwiczę przez pół godziny rano i wieczorem.", "polish"),
    @("Kierowca autobusu nosi zielone okulary przeciwsłoneczne.", "polish"),
    @("O ogrodzie zoologicznym można spotkać mnóstwo zwierząt egzotycznych.", "polish"),
    @("Mamo, nudzi mi się, co mogę zrobić?", "polish"),
    @("We have lots of free on-line games, songs, stories and activities.", "english"),
    @("This is a free app, very useful for your kids.", "english"),
    @("Are you running out of ideas on how to teach your children?", "english"),
    @("This core is theoretically enough for everyday life.", "english"),
    @("At this level students could start to move on their own.", "english"),
    @("When will a coronavirus vaccine be ready and how would it work?", "english"),
    @("Social distancing likely to go on long after the lockdown.", "english"),
    @("He thinks american music is great and the people around here are friendly.", "english"),
    @("Get a one month free trial.", "english"),
    @("Older generation must give more to the young.", "english"),
    @("Justicia prepara una ley exprés para evitar el colapso en los tribunales.", "spanish"),
    @("Se producen frutas cada vez con más azúcar.", "spanish"),
    @("Las ballenas del Atlántico tienen altos niveles de químicos procedentes del plástico.", "spanish"),
    @("No ejecutar en páginas de este sitio web", "spanish"),
    @("Barcelona es una ciudad española, capital de la comunidad autónoma de Cataluña.", "spanish"),
    @("La ciudad posee un clima mediterráneo con influencias marítimas.", "spanish"),
    @("Cada distrito funciona como un ente político con competencias propias.", "spanish"),
    @("¡Bienvenido a la única web especializada en actividades y excursiones en el Parque Nacional del Teide!", "spanish"),
    @("Atendiendo a las clases existen dos tipos de desempleo.", "spanish"),
    @("A los efectos de la base reguladora no se incluirá en la misma las horas extraordinarias.", "spanish")
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new ml_tutorial training sentences (polish / english / spanish)
# right after the existing data, starting at row 32.
$startRow = 32
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Widen column A to fit the longer sentences that were just added.
$ws.Columns.Item(1).ColumnWidth = 84.17

# Restore the view/selection state recorded for the edited sheet.
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F59").Select()
